$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Relocate the hidden "_GoBack" bookmark from the end of the
#    "Visit GitHub Website" run down in section 3 to the end of the
#    "What is Git?" heading paragraph near the top of the document,
#    and add a new run containing "?" right before it (so the
#    heading paragraph ends up with two runs: "What is Git?" + "?").
# ------------------------------------------------------------------

# Remove the existing "_GoBack" bookmark (hidden bookmarks are not
# enumerated by Bookmarks.Count/ForEach, but they can still be
# reached directly by name).
try {
    $oldGoBack = $d.Bookmarks.Item("_GoBack")
    $oldGoBack.Delete()
} catch {
    # nothing to remove
}

# Locate the "What is Git?" heading text.
$searchRange = $d.Content
$searchRange.Find.ClearFormatting()
$found = $searchRange.Find.Execute("What is Git?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertPos = $searchRange.End

    # Insert a two-character placeholder ("?X") right after the
    # existing "What is Git?" text so the new "?" is not the very
    # last character of the paragraph while we position the
    # bookmark (inserting a bookmark exactly at end-of-paragraph
    # is unreliable), then toggle a character property on the new
    # text so Word keeps it as its own run instead of merging it
    # back into the previous "What is Git?" run.
    $insertRange = $d.Range($insertPos, $insertPos)
    $insertRange.InsertAfter("?X")

    $newRunRange = $d.Range($insertPos, $insertPos + 2)
    $newRunRange.Font.Bold = 1
    $newRunRange.Font.Bold = 0

    # Add the bookmark collapsed between the real "?" and the
    # temporary "X" placeholder.
    $bookmarkPos = $insertPos + 1
    $bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)

    # Drop the temporary placeholder character; the bookmark (being
    # collapsed just before it) stays put at the paragraph end.
    $placeholderRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
    $placeholderRange.Delete()
}
